$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Product Backlog": swap backlog items (rows 2 & 3) and move the
# "Time of Record" date cell from I2 down to I3 (leaving I3 empty but with
# the same date formatting, and I2 cleared entirely).
# ---------------------------------------------------------------------------
$wsPB = $wb.Worksheets.Item("Product Backlog")

$b2 = $wsPB.Range("B2").Value()
$c2 = $wsPB.Range("C2").Value()
$e2 = $wsPB.Range("E2").Value()
$b3 = $wsPB.Range("B3").Value()
$c3 = $wsPB.Range("C3").Value()
$e3 = $wsPB.Range("E3").Value()

$wsPB.Range("B2").Value = $b3
$wsPB.Range("C2").Value = $c3
$wsPB.Range("E2").Value = $e3
$wsPB.Range("B3").Value = $b2
$wsPB.Range("C3").Value = $c2
$wsPB.Range("E3").Value = $e2

# Carry the date cell's formatting from I2 to I3 without its value, then wipe I2.
$wsPB.Range("I2").Copy() | Out-Null
$wsPB.Range("I3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsPB.Range("I2").Clear() | Out-Null

$wsPB.Range("I2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Sprint Backlog": renumber backlog IDs and update a few sprint
# status / effort values.
# ---------------------------------------------------------------------------
$wsSB = $wb.Worksheets.Item("Sprint Backlog")

$wsSB.Range("A2").Value = 1.01
$wsSB.Range("A3").Value = 1.02
$wsSB.Range("A4").Value = 1.03
$wsSB.Range("A5").Value = 1.04
$wsSB.Range("A6").Value = 1.05
$wsSB.Range("A7").Value = 1.06
$wsSB.Range("A8").Value = 1.07
$wsSB.Range("A9").Value = 1.08
$wsSB.Range("A10").Value = 1.09
$wsSB.Range("A11").Value = 1.1

$wsSB.Range("K5").Value = 5
$wsSB.Range("K8").Value = 4

$wsSB.Range("L3").Value = "done"
$wsSB.Range("L5").Value = "done"
$wsSB.Range("L6").Value = "done"
$wsSB.Range("L8").Value = "work in progress"
$wsSB.Range("L10").Value = "work in progress"

# ---------------------------------------------------------------------------
# Sheet "BurndownChart": only the selected range changed.
# ---------------------------------------------------------------------------
$wsBC = $wb.Worksheets.Item("BurndownChart")
$wsBC.Range("A3:E3").Select() | Out-Null

# Re-activate "Sprint Backlog" last so it stays the workbook's selected tab,
# matching its sheetView tabSelected="1" / the workbook's original activeTab.
$wsSB.Activate() | Out-Null
$wsSB.Range("C14").Select() | Out-Null

Write-Host "edits applied"
